# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.345.26'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '1.665.95'
$ws.Range("E3").Value = '  -2.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.64'
$ws.Range("E5").Value = '  -2.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5164'
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.008'
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06450'
$ws.Range("E8").Value = '  -2.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2565'
$ws.Range("E9").Value = '  -3.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.93'
$ws.Range("E10").Value = '  -4.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07668'
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("D12").Value = '1.672.32'
$ws.Range("E12").Value = '  -2.60%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.896.29'
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.315'
$ws.Range("E14").Value = '  -5.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5544'
$ws.Range("E15").Value = '  -3.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.46'
$ws.Range("E17").Value = '  -5.16%  '
$ws.Range("D18").Value = '26.383.26'
$ws.Range("E18").Value = '  -3.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.37'
$ws.Range("E20").Value = '  -2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.401'
$ws.Range("E21").Value = '  -5.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.11'
$ws.Range("E22").Value = '  -3.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.899'
$ws.Range("E23").Value = '  -1.41%  '
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.47'
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.754'
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("E27").Value = '  -4.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.994'
$ws.Range("E28").Value = '  -3.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.77'
$ws.Range("E29").Value = '  -3.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05253'
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.262'
$ws.Range("E31").Value = '  -2.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.380'
$ws.Range("E32").Value = '  -3.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.221'
$ws.Range("E33").Value = '  -6.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.565'
$ws.Range("E34").Value = '  -5.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.758'
$ws.Range("E35").Value = '  -4.31%  '
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9263'
$ws.Range("E37").Value = '  -2.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5729'
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("D39").Value = '1.155.23'
$ws.Range("E39").Value = '  +10.61%  '
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.007'
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8460'
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.651'
$ws.Range("E43").Value = '  -3.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.88'
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").Value = '1.806.36'
$ws.Range("E46").Value = '  -5.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4502'
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '56.11'
$ws.Range("E48").Value = '  -3.49%  '
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.928'
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("E51").Value = '  -2.71%  '
